$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells: Wins / Losses / Ties in columns AD, AE, AF of row 1
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the bold/centered/bordered style used by the other header cells (e.g. AC1)
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Fill team record values (Wins=95, Losses=67, Ties=0) for all data rows (2 through 48)
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 30).Value = 95  # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 67  # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF -> Ties
}
